$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "38.900.12"
Set-TextValue $ws.Range("E2") "  +0.29%  "
Set-TextValue $ws.Range("D3") "2.148.37"
Set-TextValue $ws.Range("E3") "  +2.64%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "228.42"
Set-TextValue $ws.Range("D6") "0.618"
Set-TextValue $ws.Range("E6") "  +0.70%  "
Set-TextValue $ws.Range("D7") "62.27"
Set-TextValue $ws.Range("E7") "  +2.63%  "
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.05%  "
Set-TextValue $ws.Range("D9") "0.391"
Set-TextValue $ws.Range("E9") "  +1.61%  "
Set-TextValue $ws.Range("E10") "  +1.17%  "
Set-TextValue $ws.Range("E11") "  -0.21%  "
Set-TextValue $ws.Range("D12") "15.98"
Set-TextValue $ws.Range("E12") "  +6.75%  "
Set-TextValue $ws.Range("D13") "2.465.76"
Set-TextValue $ws.Range("E13") "  +2.60%  "
Set-TextValue $ws.Range("D14") "22.17"
Set-TextValue $ws.Range("E14") "  +1.06%  "
Set-TextValue $ws.Range("D15") "0.811"
Set-TextValue $ws.Range("E15") "  +1.85%  "
Set-TextValue $ws.Range("E16") "  +0.91%  "
Set-TextValue $ws.Range("D17") "2.144.03"
Set-TextValue $ws.Range("E17") "  +2.36%  "
Set-TextValue $ws.Range("D18") "38.906.78"
Set-TextValue $ws.Range("E18") "  +0.41%  "
Set-TextValue $ws.Range("B19") "Uniswap"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "6.14"
Set-TextValue $ws.Range("E19") "  +2.00%  "
Set-TextValue $ws.Range("B20") "Litecoin"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D20") "71.85"
Set-TextValue $ws.Range("E20") "  +0.18%  "
Set-TextValue $ws.Range("D21") "0.0₃0847"
Set-TextValue $ws.Range("E21") "  +1.30%  "
Set-TextValue $ws.Range("D22") "227.61"
Set-TextValue $ws.Range("E22") "  +0.58%  "
Set-TextValue $ws.Range("E23") "  -0.03%  "
Set-TextValue $ws.Range("D24") "2.41"
Set-TextValue $ws.Range("E24") "  -0.50%  "
Set-TextValue $ws.Range("E25") "  -0.33%  "
Set-TextValue $ws.Range("D26") "9.71"
Set-TextValue $ws.Range("E26") "  +2.95%  "
Set-TextValue $ws.Range("D27") "170.42"
Set-TextValue $ws.Range("E27") "  -0.01%  "
Set-TextValue $ws.Range("E28") "  -0.29%  "
Set-TextValue $ws.Range("D29") "19.60"
Set-TextValue $ws.Range("E29") "  +2.30%  "
Set-TextValue $ws.Range("D30") "1.40"
Set-TextValue $ws.Range("E30") "  -3.20%  "
Set-TextValue $ws.Range("E31") "  +9.69%  "
Set-TextValue $ws.Range("E32") "  +0.83%  "
Set-TextValue $ws.Range("D34") "4.80"
Set-TextValue $ws.Range("E34") "  +2.00%  "
Set-TextValue $ws.Range("E35") "  +11.31%  "
Set-TextValue $ws.Range("E36") "  +0.40%  "
Set-TextValue $ws.Range("D37") "2.42"
Set-TextValue $ws.Range("E37") "  +1.44%  "
Set-TextValue $ws.Range("D38") "3.54"
Set-TextValue $ws.Range("E38") "  +0.65%  "
Set-TextValue $ws.Range("D39") "1.00"
Set-TextValue $ws.Range("E39") "  -0.01%  "
Set-TextValue $ws.Range("D40") "18.17"
Set-TextValue $ws.Range("E40") "  -0.13%  "
Set-TextValue $ws.Range("D42") "102.78"
Set-TextValue $ws.Range("E42") "  +1.72%  "
Set-TextValue $ws.Range("D43") "1.536.10"
Set-TextValue $ws.Range("E43") "  -0.19%  "
Set-TextValue $ws.Range("E44") "  +6.56%  "
Set-TextValue $ws.Range("B45") "ARBITRUM"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D45") "1.10"
Set-TextValue $ws.Range("E45") "  +6.84%  "
Set-TextValue $ws.Range("B46") "FraxShare"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D46") "7.83"
Set-TextValue $ws.Range("E46") "  +2.38%  "
Set-TextValue $ws.Range("B47") "HuobiToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D47") "2.80"
Set-TextValue $ws.Range("E47") "  -0.51%  "
Set-TextValue $ws.Range("D48") "0.0916"
Set-TextValue $ws.Range("E48") "  -0.88%  "
Set-TextValue $ws.Range("E49") "  +1.67%  "
Set-TextValue $ws.Range("D50") "2.99"
Set-TextValue $ws.Range("E50") "  +0.55%  "
Set-TextValue $ws.Range("D51") "2.350.84"
Set-TextValue $ws.Range("E51") "  +2.70%  "
